# Weekly update: insert the newest Coliflor price record at the top of the
# date-ordered data block (row 555), pushing all existing records (rows
# 555-610) down by one row (to 556-611).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 555; everything below (old 555..610) shifts to 556..611.
$ws.Rows.Item(555).Insert()

# Populate the newly inserted row with this week's data point.
$ws.Range("A555").Value = 4
$ws.Range("B555").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C555").Value = "Los Lagos"
$ws.Range("D555").Value = 45212
$ws.Range("E555").Value = 10
$ws.Range("F555").Value = 100112008
$ws.Range("G555").Value = "Coliflor"
$ws.Range("H555").Value = "Sin especificar"
$ws.Range("I555").Value = "Primera"
$ws.Range("J555").Value = 1500
$ws.Range("K555").Value = 1500
$ws.Range("L555").Value = 1500
$ws.Range("M555").Value = 1500
$ws.Range("N555").Value = "$/unidad"
$ws.Range("O555").Value = "Región Metropolitana"
$ws.Range("P555").Value = 1500
$ws.Range("Q555").Value = 1
$ws.Range("R555").Value = "Hortaliza"
